$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet ---
$ws.Name = "resultados"

# --- Halve the "Primera_Centrada" (column D) values ---
$ws.Range("D2").Value = -3.7187999999999999
$ws.Range("D3").Value = 7.25
$ws.Range("D4").Value = 4.8125
$ws.Range("D5").Value = 2.75
$ws.Range("D6").Value = 1.0625
$ws.Range("D7").Value = -0.25
$ws.Range("D8").Value = -1.1875
$ws.Range("D9").Value = -1.75
$ws.Range("D10").Value = -1.9375
$ws.Range("D11").Value = -1.75
$ws.Range("D12").Value = -1.1875
$ws.Range("D13").Value = -0.25
$ws.Range("D14").Value = 1.0625
$ws.Range("D15").Value = 2.75
$ws.Range("D16").Value = 4.8125
$ws.Range("D17").Value = 7.25
$ws.Range("D18").Value = -3.7187999999999999

# --- Fix the typo in the second chart's title: "Sugunda" -> "Segunda" ---
$co2Title = $ws.ChartObjects(2)
$co2Title.Chart.ChartTitle.Text = "Diferencias finitas (Segunda derivada)"

# --- Reposition the two charts to match the updated layout ---
$co1 = $ws.ChartObjects(1)
$co1.Left = 467.5
$co1.Top = 0.0
$co1.Width = 340.1875
$co1.Height = 216.00001

$co2 = $ws.ChartObjects(2)
$co2.Left = 467.5
$co2.Top = 230.1
$co2.Width = 340.1875
$co2.Height = 216.0

Write-Output "done"
